# aalatih.xlsx — update the leak-pressure readings for rows 134-143 (column C)
# from 21 to 21.5, and move the sheet's view/selection down to D137
# (with the window scrolled so row 127 is at the top), matching the
# state the workbook was left in after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column C (Titik_1_PSI) for rows 134 through 143: 21 -> 21.5
for ($r = 134; $r -le 143; $r++) {
    $ws.Cells.Item($r, 3).Value = 21.5
}

# Make Sheet1 the active sheet/tab and reproduce the final selection +
# scroll position (top-left visible cell A127, active cell D137).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 127
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D137").Select()
